# Commit: Tue, Jun 16, 2020  1:04:52 AM
#
# 1) Table on slide 5 switches to a different built-in table style.
# 2) The presentation's theme (color scheme) changes from the "Integral" /
#    "Red Violet" design to the plain "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{3F464DBC-F52F-4E08-B50F-95639EB1A133}")
    }
}

# --- 2. Switch the theme colours back to the default Office palette -------
$tcs = $slide5.ThemeColorScheme
$tcs.Item(1).RGB  = 0            # dk1      000000
$tcs.Item(2).RGB  = 16777215     # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388      # dk2      44546A
$tcs.Item(4).RGB  = 15132391     # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939     # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501      # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845     # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407        # accent4  FFC000
$tcs.Item(9).RGB  = 12874308     # accent5  4472C4
$tcs.Item(10).RGB = 4697456      # accent6  70AD47
$tcs.Item(11).RGB = 12673797     # hlink    0563C1
$tcs.Item(12).RGB = 7491477      # folHlink 954F72
